$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.4341637010676156
$ws1.Range("C2").Value = 0.07848837209302326
$ws1.Range("D2").Value = 0.9642857142857143
$ws1.Range("E2").Value = 0.1451612903225807
$ws1.Range("F2").Value = 0.2960526315789473
$ws1.Range("G2").Value = 0.6724137931034483
$ws1.Range("H2").Value = 0.7934724451578384
$ws1.Range("I2").Value = 27
$ws1.Range("J2").Value = 317
$ws1.Range("K2").Value = 217
$ws1.Range("L2").Value = 1

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9954128440366973
$ws2.Range("C2").Value = 0.4063670411985019
$ws2.Range("D2").Value = 0.5771276595744681

$ws2.Range("B3").Value = 0.07848837209302326
$ws2.Range("C3").Value = 0.9642857142857143
$ws2.Range("D3").Value = 0.1451612903225807

$ws2.Range("B4").Value = 0.4341637010676156
$ws2.Range("C4").Value = 0.4341637010676156
$ws2.Range("D4").Value = 0.4341637010676156
$ws2.Range("E4").Value = 0.4341637010676156

$ws2.Range("B5").Value = 0.5369506080648603
$ws2.Range("C5").Value = 0.6853263777421081
$ws2.Range("D5").Value = 0.3611444749485244

$ws2.Range("B6").Value = 0.9497297742601442
$ws2.Range("C6").Value = 0.4341637010676156
$ws2.Range("D6").Value = 0.555606203455157

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 217
$ws3.Range("C2").Value = 317
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 27
